$d = $word.ActiveDocument

# Helper: find the paragraph index whose trimmed text equals $target exactly
function Find-ParaIndex($doc, $target) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        # strip trailing paragraph mark / cell mark characters
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $target) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1. Move the hidden "_GoBack" bookmark from the end of the "Clean up
#    tests." paragraph to a brand-new, otherwise-empty paragraph that
#    is inserted immediately before "Light Windows API wrapper.".
# ------------------------------------------------------------------
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Insert the new "Use BOOST_TEST_DONT_PRINT_LOG_VALUE where
#    appropriate." bullet right after "Add concept checks in unit
#    tests for constructability, moveability, destructibility, etc."
# ------------------------------------------------------------------
$idx = Find-ParaIndex $d "Add concept checks in unit tests for constructability, moveability, destructibility, etc."
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($idx + 1)
$newPara.Range.Text = "Use BOOST_TEST_DONT_PRINT_LOG_VALUE where appropriate."

# ------------------------------------------------------------------
# 3. Remove the "Stack trace in exception object." bullet entirely.
# ------------------------------------------------------------------
$idx = Find-ParaIndex $d "Stack trace in exception object."
$p = $d.Paragraphs.Item($idx)
$p.Range.Delete()

# ------------------------------------------------------------------
# 4. Insert a new (empty) top-level bullet paragraph right before
#    "Light Windows API wrapper." and re-home the "_GoBack" bookmark
#    there (collapsed, no visible text), matching the target XML:
#      <w:p>...<w:bookmarkStart .../><w:bookmarkEnd .../></w:p>
# ------------------------------------------------------------------
$idx = Find-ParaIndex $d "Light Windows API wrapper."
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphBefore()
$bmPara = $d.Paragraphs.Item($idx)

# Use a temporary placeholder character so Bookmarks.Add gets a
# non-degenerate range, then delete the character - the bookmark
# start/end tags collapse in place and stay in this paragraph.
$bmPara.Range.Text = "X"
$bmRange = $d.Range($bmPara.Range.Start, $bmPara.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$clearRange = $d.Range($bmPara.Range.Start, $bmPara.Range.Start + 1)
$clearRange.Text = ""

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
